# Regenerate column G ("K") values for each data row (rows 2-32).
# This mirrors the upstream save_data regeneration which now uses K
# (strike count, i.e. s_vals) instead of the old Strike# values, after
# recalculating std/mean for the underlying distribution.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 5
    3  = 6
    4  = 2
    5  = 7
    6  = 2
    7  = 2
    8  = 6
    9  = 4
    10 = 0
    11 = 2
    12 = 2
    13 = 7
    14 = 2
    15 = 7
    16 = 1
    17 = 3
    18 = 3
    19 = 1
    20 = 6
    21 = 2
    22 = 2
    23 = 3
    24 = 4
    25 = 4
    26 = 4
    27 = 4
    28 = 4
    29 = 1
    30 = 3
    31 = 0
    32 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
